$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I and J, matching style/format of existing
# header H1 (bold font, thin border, centered alignment)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-46 for columns I (I0) and J (IF)
$data = @(
    @(2, 8, 9),
    @(3, 6, 7),
    @(4, 6, 6),
    @(5, 1, 1),
    @(6, 9, 9),
    @(7, 10, 10),
    @(8, 9, 9),
    @(9, 1, 1),
    @(10, 6, 6),
    @(11, 7, 7),
    @(12, 4, 5),
    @(13, 5, 5),
    @(14, 3, 4),
    @(15, 6, 6),
    @(16, 5, 5),
    @(17, 7, 7),
    @(18, 7, 7),
    @(19, 11, 11),
    @(20, 4, 5),
    @(21, 9, 9),
    @(22, 7, 7),
    @(23, 5, 5),
    @(24, 10, 10),
    @(25, 6, 7),
    @(26, 5, 6),
    @(27, 7, 7),
    @(28, 8, 8),
    @(29, 10, 10),
    @(30, 1, 2),
    @(31, 8, 8),
    @(32, 7, 8),
    @(33, 1, 2),
    @(34, 10, 10),
    @(35, 1, 2),
    @(36, 8, 8),
    @(37, 5, 5),
    @(38, 6, 6),
    @(39, 8, 9),
    @(40, 6, 6),
    @(41, 7, 7),
    @(42, 7, 8),
    @(43, 3, 3),
    @(44, 7, 7),
    @(45, 5, 5),
    @(46, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
